# Auto-generated Excel COM-interop script to apply scheduled-runner value updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.

$wb = $excel.ActiveWorkbook

function Set-CellValue($ws, $addr, $value) {
    $ws.Range($addr).Value = $value
}

function Clear-CellValue($ws, $addr) {
    $ws.Range($addr).ClearContents()
}


# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 51
Set-CellValue $ws "H51" 2996.8413   # was 2998.4033
Set-CellValue $ws "J51" 2989   # was 3000.125
Set-CellValue $ws "L51" 2989   # was 3000.125
Set-CellValue $ws "N51" -3957   # was -3968.125
# Row 135
Set-CellValue $ws "H135" 2427.5334   # was 2593.9312
Set-CellValue $ws "I135" 2173.16   # was 2292.6086
Set-CellValue $ws "J135" 3699.4   # was 3749
Set-CellValue $ws "K135" 19558.44   # was 20633.4774
Set-CellValue $ws "L135" 33294.6   # was 33741
Set-CellValue $ws "M135" -17023.44   # was -18098.4774
Set-CellValue $ws "N135" -38364.6   # was -38811
# Row 140
Set-CellValue $ws "H140" 0   # was 100000
Set-CellValue $ws "J140" 0   # was 100000
Set-CellValue $ws "L140" 0   # was 100000
Clear-CellValue $ws "N140"   # was -110360

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 4
Set-CellValue $ws "H4" 0   # was 196
Set-CellValue $ws "I4" 0   # was 227.5
Set-CellValue $ws "J4" 0   # was 70
Set-CellValue $ws "K4" 0   # was 227.5
Set-CellValue $ws "L4" 0   # was 70
Clear-CellValue $ws "M4"   # was -111.5
Clear-CellValue $ws "N4"   # was -302
# Row 32
Set-CellValue $ws "H32" 1050834   # was 1050836.1
Set-CellValue $ws "I32" 1236446.5   # was 1236449
Set-CellValue $ws "K32" 1236446.5   # was 1236449
Set-CellValue $ws "M32" -1236159.5   # was -1236162
# Row 45
Set-CellValue $ws "H45" 2924.95   # was 3213.2778
Set-CellValue $ws "I45" 2684.0625   # was 2937.4666
Set-CellValue $ws "J45" 3888.5   # was 4592.3335
Set-CellValue $ws "K45" 2684.0625   # was 2937.4666
Set-CellValue $ws "L45" 3888.5   # was 4592.3335
Set-CellValue $ws "M45" -2307.0625   # was -2560.4666
Set-CellValue $ws "N45" -4642.5   # was -5346.3335
# Row 74
Set-CellValue $ws "H74" 1073747.8   # was 997083.75
Set-CellValue $ws "I74" 1295019.5   # was 1210590.5
Set-CellValue $ws "J74" 16560.777   # was 14952.7
Set-CellValue $ws "K74" 1295019.5   # was 1210590.5
Set-CellValue $ws "L74" 16560.777   # was 14952.7
Set-CellValue $ws "M74" -1294145.5   # was -1209716.5
Set-CellValue $ws "N74" -18308.777   # was -16700.7
# Row 77
Set-CellValue $ws "H77" 1073747.8   # was 997083.75
Set-CellValue $ws "I77" 1295019.5   # was 1210590.5
Set-CellValue $ws "J77" 16560.777   # was 14952.7
Set-CellValue $ws "K77" 6475097.5   # was 6052952.5
Set-CellValue $ws "L77" 82803.88499999999   # was 74763.5
Set-CellValue $ws "M77" -6470729.5   # was -6048584.5
Set-CellValue $ws "N77" -91539.88499999999   # was -83499.5
# Row 132
Set-CellValue $ws "H132" 611918.75   # was 678017.4399999999
Set-CellValue $ws "I132" 659975.8   # was 737549.9
Set-CellValue $ws "J132" 3196.3333   # was 3316.6667
Set-CellValue $ws "K132" 1979927.4   # was 2212649.7
Set-CellValue $ws "L132" 9588.999899999999   # was 9950.000100000001
Set-CellValue $ws "M132" -1977397.4   # was -2210119.7
Set-CellValue $ws "N132" -14648.9999   # was -15010.0001

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
Set-CellValue $ws "H20" 30037.703   # was 30866.916
Set-CellValue $ws "I20" 37588.242   # was 40332.406
Set-CellValue $ws "J20" 2667   # was 2470.4443
Set-CellValue $ws "K20" 37588.242   # was 40332.406
Set-CellValue $ws "L20" 2667   # was 2470.4443
Set-CellValue $ws "M20" -37341.242   # was -40085.406
Set-CellValue $ws "N20" -3161   # was -2964.4443
# Row 94
Set-CellValue $ws "H94" 13133.333   # was 4510
Set-CellValue $ws "I94" 12200   # was 4022.4285
Set-CellValue $ws "J94" 15000   # was 5647.6665
Set-CellValue $ws "K94" 12200   # was 4022.4285
Set-CellValue $ws "L94" 15000   # was 5647.6665
Set-CellValue $ws "M94" -11749   # was -3571.4285
Set-CellValue $ws "N94" -15902   # was -6549.6665
# Row 107
Set-CellValue $ws "H107" 1501.9286   # was 1531.1482
Set-CellValue $ws "J107" 1356.5   # was 2000
Set-CellValue $ws "L107" 1356.5   # was 2000
Set-CellValue $ws "N107" -5196.5   # was -5840

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
Set-CellValue $ws "H22" 1583.7   # was 1518
Set-CellValue $ws "I22" 168.4   # was 178.8
Set-CellValue $ws "J22" 2999   # was 2735.4546
Set-CellValue $ws "K22" 168.4   # was 178.8
Set-CellValue $ws "L22" 2999   # was 2735.4546
Set-CellValue $ws "M22" 181.6   # was 171.2
Set-CellValue $ws "N22" -3699   # was -3435.4546
# Row 87
Set-CellValue $ws "H87" 0   # was 60000
Set-CellValue $ws "J87" 0   # was 60000
Set-CellValue $ws "L87" 0   # was 60000
Clear-CellValue $ws "N87"   # was -62372
# Row 90
Set-CellValue $ws "H90" 0   # was 60000
Set-CellValue $ws "J90" 0   # was 60000
Set-CellValue $ws "L90" 0   # was 180000
Clear-CellValue $ws "N90"   # was -191856
# Row 92
Set-CellValue $ws "H92" 32500   # was 38333.332
Set-CellValue $ws "J92" 32500   # was 38333.332
Set-CellValue $ws "L92" 32500   # was 38333.332
Set-CellValue $ws "N92" -37492   # was -43325.332
# Row 107
Set-CellValue $ws "H107" 1763.3334   # was 2269.0908
Set-CellValue $ws "I107" 1316   # was 1496.25
Set-CellValue $ws "J107" 4000   # was 4330
Set-CellValue $ws "K107" 1316   # was 1496.25
Set-CellValue $ws "L107" 4000   # was 4330
Set-CellValue $ws "M107" 604   # was 423.75
Set-CellValue $ws "N107" -7840   # was -8170
# Row 132
Set-CellValue $ws "H132" 1849.68   # was 1911.1666
Set-CellValue $ws "I132" 1849.68   # was 1911.1666
Set-CellValue $ws "K132" 5549.04   # was 5733.4998
Set-CellValue $ws "M132" -3019.04   # was -3203.4998
# Row 134
Set-CellValue $ws "H134" 2750.173   # was 2824.4082
Set-CellValue $ws "I134" 2599.2058   # was 2680.4375
Set-CellValue $ws "J134" 3035.3333   # was 3095.4119
Set-CellValue $ws "K134" 7797.617400000001   # was 8041.3125
Set-CellValue $ws "L134" 9105.999899999999   # was 9286.235700000001
Set-CellValue $ws "M134" -5262.617400000001   # was -5506.3125
Set-CellValue $ws "N134" -14175.9999   # was -14356.2357

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
Set-CellValue $ws "H4" 3606149.8   # was 4124057.8
Set-CellValue $ws "I4" 4588636   # was 5183578.5
Set-CellValue $ws "K4" 13765908   # was 15550735.5
Set-CellValue $ws "M4" -13765796   # was -15550623.5
# Row 17
Set-CellValue $ws "H17" 1180   # was 1018.3333
Set-CellValue $ws "I17" 45   # was 122.5
Set-CellValue $ws "J17" 1463.75   # was 1466.25
Set-CellValue $ws "K17" 135   # was 367.5
Set-CellValue $ws "L17" 4391.25   # was 4398.75
Set-CellValue $ws "M17" 34   # was -198.5
Set-CellValue $ws "N17" -4729.25   # was -4736.75
# Row 40
Set-CellValue $ws "H40" 44.904762   # was 52.42857
Set-CellValue $ws "I40" 33.9   # was 60
Set-CellValue $ws "J40" 54.909092   # was 45.545456
Set-CellValue $ws "K40" 135.6   # was 240
Set-CellValue $ws "L40" 219.636368   # was 182.181824
Set-CellValue $ws "M40" -66.59999999999999   # was -171
Set-CellValue $ws "N40" -357.636368   # was -320.181824
# Row 122
Set-CellValue $ws "H122" 1613023.4   # was 1075371.9
Set-CellValue $ws "I122" 2016259.9   # was 1152180.6
Set-CellValue $ws "J122" 77.5   # was 50
Set-CellValue $ws "K122" 18146339.1   # was 10369625.4
Set-CellValue $ws "L122" 697.5   # was 450
Set-CellValue $ws "M122" -18143889.1   # was -10367175.4
Set-CellValue $ws "N122" -5597.5   # was -5350

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
Set-CellValue $ws "H2" 3125109.8   # was 3125109.5
Set-CellValue $ws "I2" 4347871.5   # was 4166710.5
Set-CellValue $ws "J2" 273.22223   # was 307.125
Set-CellValue $ws "K2" 4347871.5   # was 4166710.5
Set-CellValue $ws "L2" 273.22223   # was 307.125
Set-CellValue $ws "M2" -4347758.5   # was -4166597.5
Set-CellValue $ws "N2" -499.22223   # was -533.125
# Row 132
Set-CellValue $ws "H132" 17321.121   # was 17768.875
Set-CellValue $ws "I132" 18164.42   # was 18196.902
Set-CellValue $ws "J132" 4250   # was 4500
Set-CellValue $ws "K132" 54493.25999999999   # was 54590.70599999999
Set-CellValue $ws "L132" 12750   # was 13500
Set-CellValue $ws "M132" -51963.25999999999   # was -52060.70599999999
Set-CellValue $ws "N132" -17810   # was -18560

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
Set-CellValue $ws "H16" 2406.0557   # was 2495.2354
Set-CellValue $ws "I16" 593.6923   # was 586.1539
Set-CellValue $ws "J16" 7118.2   # was 8699.75
Set-CellValue $ws "K16" 593.6923   # was 586.1539
Set-CellValue $ws "L16" 7118.2   # was 8699.75
Set-CellValue $ws "M16" -423.6923   # was -416.1539
Set-CellValue $ws "N16" -7458.2   # was -9039.75
# Row 38
Set-CellValue $ws "H38" 15000   # was 0
Set-CellValue $ws "J38" 15000   # was 0
Set-CellValue $ws "L38" 15000   # was 0
Set-CellValue $ws "N38" -15820   # newly added (was empty)
# Row 55
Set-CellValue $ws "H55" 1119.825   # was 1250.0571
Set-CellValue $ws "I55" 1031.1052   # was 1252.8667
Set-CellValue $ws "J55" 1200.0952   # was 1247.95
Set-CellValue $ws "K55" 1031.1052   # was 1252.8667
Set-CellValue $ws "L55" 1200.0952   # was 1247.95
Set-CellValue $ws "M55" -858.1052   # was -1079.8667
Set-CellValue $ws "N55" -1546.0952   # was -1593.95
# Row 68
Set-CellValue $ws "H68" 16840.45   # was 100926.57
Set-CellValue $ws "I68" 13390.363   # was 152885.88
Set-CellValue $ws "J68" 27683.572   # was 31647.5
Set-CellValue $ws "K68" 13390.363   # was 152885.88
Set-CellValue $ws "L68" 27683.572   # was 31647.5
Set-CellValue $ws "M68" -12641.363   # was -152136.88
Set-CellValue $ws "N68" -29181.572   # was -33145.5
# Row 71
Set-CellValue $ws "H71" 16840.45   # was 100926.57
Set-CellValue $ws "I71" 13390.363   # was 152885.88
Set-CellValue $ws "J71" 27683.572   # was 31647.5
Set-CellValue $ws "K71" 66951.815   # was 764429.4
Set-CellValue $ws "L71" 138417.86   # was 158237.5
Set-CellValue $ws "M71" -63207.815   # was -760685.4
Set-CellValue $ws "N71" -145905.86   # was -165725.5
# Row 133
Set-CellValue $ws "H133" 89317.5   # was 89316
Set-CellValue $ws "J133" 89317.5   # was 89316
Set-CellValue $ws "L133" 89317.5   # was 89316
Set-CellValue $ws "N133" -94377.5   # was -94376

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 45
Set-CellValue $ws "H45" 16523.143   # was 16880.428
Set-CellValue $ws "I45" 15964.333   # was 16631.334
Set-CellValue $ws "J45" 16942.25   # was 17067.25
Set-CellValue $ws "K45" 15964.333   # was 16631.334
Set-CellValue $ws "L45" 16942.25   # was 17067.25
Set-CellValue $ws "M45" -15473.333   # was -16140.334
Set-CellValue $ws "N45" -17924.25   # was -18049.25
# Row 54
Set-CellValue $ws "H54" 25000   # was 19333.334
Set-CellValue $ws "I54" 25000   # was 19333.334
Set-CellValue $ws "K54" 25000   # was 19333.334
Set-CellValue $ws "M54" -24480   # was -18813.334
# Row 74
Set-CellValue $ws "H74" 10013190   # was 11124711
Set-CellValue $ws "J74" 14655.444   # was 15300.25
Set-CellValue $ws "L74" 14655.444   # was 15300.25
Set-CellValue $ws "N74" -16527.444   # was -17172.25
# Row 77
Set-CellValue $ws "H77" 10013190   # was 11124711
Set-CellValue $ws "J77" 14655.444   # was 15300.25
Set-CellValue $ws "L77" 43966.33199999999   # was 45900.75
Set-CellValue $ws "N77" -53326.33199999999   # was -55260.75
# Row 112
Set-CellValue $ws "H112" 120000   # was 56795.668
Set-CellValue $ws "I112" 0   # was 50000
Set-CellValue $ws "J112" 120000   # was 60193.5
Set-CellValue $ws "K112" 0   # was 50000
Set-CellValue $ws "L112" 120000   # was 60193.5
Clear-CellValue $ws "M112"   # was -48523
Set-CellValue $ws "N112" -122954   # was -63147.5
